$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.5448
$ws.Range("B3").Value = 5.965699999999988
$ws.Range("C5").Value = -14.32120000000001
$ws.Range("D5").Value = -8.451999999999998
$ws.Range("D9").Value = -7.2811
$ws.Range("D11").Value = -8.336500000000008
$ws.Range("B14").Value = 9.057399999999999
$ws.Range("B16").Value = 9.811499999999999
$ws.Range("C16").Value = -11.48329999999999
$ws.Range("D17").Value = -7.1882
$ws.Range("B21").Value = 5.739699999999991
$ws.Range("D21").Value = -7.862800000000002
$ws.Range("B23").Value = 5.492500000000001
$ws.Range("B25").Value = 5.693099999999994
